$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[-, 'MEC-2NA-C.L.P.', 'MEC-2NB-Coman. Hidraulicos', -]"
$ws.Range("C18").Value = "[-, 'ELM-2NA-Eletropneumática', -, 'MEC-2NB-C. L. P.']"
$ws.Range("E18").Value = "[-, 'MEC-2NB-C. L. P.', 'MEC-2NB-Coman. Hidraulicos', -]"
$ws.Range("F18").Value = "[-, 'MEC-2NA-C.L.P.', 'MEC-2NA-C. Hidráulica', -]"

$ws.Range("B19").Value = "['ELM-2NA-Eletro', 'MEC-2NA-C.pneumática', 'ELM-2NA-Eletropneumática', -]"
$ws.Range("C19").Value = "ELM-1NA-Máquinas Térmicas e de Fluxo"
$ws.Range("D19").Value = "[-, -, -, 'MEC-2NA-C.pneumática']"
$ws.Range("E19").Value = "[-, 'MEC-2NB-C. L. P.', 'MEC-2NB-Coman. Hidraulicos', -]"
$ws.Range("F19").Value = "[-, 'MEC-2NA-C.L.P.', 'MEC-2NA-C. Hidráulica', 'MEC-2NA-C.pneumática']"

$ws.Range("B20").Value = "['ELM-2NA-Eletro', 'MEC-2NA-C.pneumática', -, -]"
$ws.Range("C20").Value = "['ELM-2NA-Eletro', 'ELM-2NA-Eletropneumática', -, -]"
$ws.Range("E20").Value = "ELM-1NA-Máquinas Térmicas e de Fluxo"
$ws.Range("F20").Value = "[-, 'MEC-2NA-C.L.P.', 'MEC-2NA-C. Hidráulica', -]"

$ws.Range("B21").Value = "['ELM-2NA-Eletro', -, -, 'ELM-2NA-Eletropneumática']"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "['MEC-2NB-C. L. P.', -, 'MEC-2NB-Coman. Hidraulicos', -]"
$ws.Range("F21").Value = "[-, -, 'MEC-2NA-C. Hidráulica', -]"
